# Collapse previously-split a:r runs (one word/space per run) into a
# single run per paragraph, for the specific text frames called out in
# the diff. PowerPoint's TextRange.Text setter is a no-op (keeps the
# existing runs untouched) when the new value already equals the
# current concatenated text, so each target is nudged through a throw-
# away value first to force the rewrite into one run, then set to the
# real final text.

$p = $ppt.ActivePresentation

function Set-MergedText($textRange, [string]$finalText) {
    $textRange.Text = "zzz_tmp_placeholder_zzz"
    $textRange.Text = $finalText
}

# Slide 6 ("Blank" layout slide) notes: "Blank slides can have background images."
$notes = $p.Slides.Item(6).NotesPage
Set-MergedText $notes.Shapes.Item(2).TextFrame.TextRange "Blank slides can have background images."

# Slide 1: Section Header title
Set-MergedText $p.Slides.Item(1).Shapes.Item(1).TextFrame.TextRange "Section Header (with background image)"

# Slide 2: "Slide 1" title
Set-MergedText $p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange "Slide 1"

# Slide 3: "Slide 2" title
Set-MergedText $p.Slides.Item(3).Shapes.Item(1).TextFrame.TextRange "Slide 2"

# Slide 4: "Slide 3" title
Set-MergedText $p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange "Slide 3"

# Slide 5: "Slide 4" title
Set-MergedText $p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange "Slide 4"

# Slide 5: "An image" textbox caption
Set-MergedText $p.Slides.Item(5).Shapes.Item(4).TextFrame.TextRange "An image"
